$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking price strings
# (e.g. "1.002", "23.10") are preserved exactly as text, matching
# the original inlineStr cell contents instead of being parsed as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '24.799.67'
$ws.Range("E2").Value = '  +0.77%  '
$ws.Range("D3").Value = '1.706.67'
$ws.Range("E3").Value = '  +0.98%  '
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").Value = '314.87'
$ws.Range("E5").Value = '  +0.41%  '
$ws.Range("E6").Value = '  -0.21%  '
$ws.Range("D7").Value = '0.4009'
$ws.Range("E7").Value = '  +2.76%  '
$ws.Range("D8").Value = '0.4040'
$ws.Range("E8").Value = '  +0.31%  '
$ws.Range("E9").Value = '  -0.20%  '
$ws.Range("E10").Value = '  -1.81%  '
$ws.Range("D11").Value = '53.64'
$ws.Range("E11").Value = '  +1.96%  '
$ws.Range("D12").Value = '0.08810'
$ws.Range("E12").Value = '  +0.69%  '
$ws.Range("D13").Value = '26.29'
$ws.Range("E13").Value = '  +6.16%  '
$ws.Range("D14").Value = '7.517'
$ws.Range("E14").Value = '  -0.67%  '
$ws.Range("D15").Value = '7.999'
$ws.Range("E15").Value = '  +0.57%  '
$ws.Range("D16").Value = '0.00001342'
$ws.Range("E16").Value = '  -0.66%  '
$ws.Range("D17").Value = '1.732.64'
$ws.Range("E17").Value = '  +3.04%  '
$ws.Range("D18").Value = '95.64'
$ws.Range("E18").Value = '  -2.79%  '
$ws.Range("D19").Value = '0.07170'
$ws.Range("E19").Value = '  +0.87%  '
$ws.Range("E20").Value = '  +5.53%  '
$ws.Range("D21").Value = '7.297'
$ws.Range("E21").Value = '  +0.30%  '
$ws.Range("E22").Value = '  -0.15%  '
$ws.Range("D23").Value = '14.45'
$ws.Range("E23").Value = '  +1.56%  '
$ws.Range("D24").Value = '24.789.11'
$ws.Range("E24").Value = '  +0.74%  '
$ws.Range("D25").Value = '2.353'
$ws.Range("E25").Value = '  -0.03%  '
$ws.Range("E26").Value = '  -3.41%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '23.11'
$ws.Range("E27").Value = '  +1.47%  '
$ws.Range("B28").Value = 'HuobiToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D28").Value = '6.187'
$ws.Range("E28").Value = '  +18.51%  '
$ws.Range("D29").Value = '161.81'
$ws.Range("E29").Value = '  +0.10%  '
$ws.Range("D30").Value = '144.06'
$ws.Range("E30").Value = '  +5.53%  '
$ws.Range("D31").Value = '8.228'
$ws.Range("E31").Value = '  -4.74%  '
$ws.Range("D33").Value = '1.913.52'
$ws.Range("E33").Value = '  +2.54%  '
$ws.Range("D34").Value = '0.08652'
$ws.Range("E34").Value = '  -1.60%  '
$ws.Range("E35").Value = '  +10.24%  '
$ws.Range("D36").Value = '7.303'
$ws.Range("E36").Value = '  -0.94%  '
$ws.Range("D37").Value = '1.031'
$ws.Range("E37").Value = '  -0.73%  '
$ws.Range("D38").Value = '0.2855'
$ws.Range("E38").Value = '  +4.74%  '
$ws.Range("D39").Value = '0.8410'
$ws.Range("E39").Value = '  +7.78%  '
$ws.Range("D40").Value = '0.09458'
$ws.Range("E40").Value = '  +3.53%  '
$ws.Range("D41").Value = '10.69'
$ws.Range("E41").Value = '  -0.66%  '
$ws.Range("E42").Value = '  +0.43%  '
$ws.Range("D43").Value = '1.479'
$ws.Range("E43").Value = '  +1.67%  '
$ws.Range("D44").Value = '17.51'
$ws.Range("E44").Value = '  +5.42%  '
$ws.Range("D45").Value = '2.724'
$ws.Range("E45").Value = '  +5.42%  '
$ws.Range("E46").Value = '  +3.31%  '
$ws.Range("D47").Value = '4.217'
$ws.Range("E47").Value = '  +0.66%  '
$ws.Range("D48").Value = '1.373'
$ws.Range("E48").Value = '  +3.47%  '
$ws.Range("E49").Value = '  -0.18%  '
$ws.Range("D50").Value = '140.54'
$ws.Range("E50").Value = '  +2.11%  '
$ws.Range("D51").Value = '0.08381'
$ws.Range("E51").Value = '  +5.34%  '
